# Leetcode everyday workbook update:
#  "No528. Random Pick with Weight finished and No81 reviewed"
#
# 1. Update the "knowledge point" note on row 49 (E49) from the generic
#    "数学" (math) tag to the more specific "数学，随机选取".
# 2. Append row 51 for "No528. Random Pick with Weight" (same visual
#    style as row 49).
# 3. Append row 52 for "81. Search in Rotated Sorted Array II" (same
#    visual style as row 50).
# 4. Add hyperlinks on the new C51/C52 cells.
# 5. Leave the selection on I53, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 51 -----------------------------------------------------------
# Copy the formatting block of row 49 first (same look: fill/border/etc).
$ws.Range("A49:H49").Copy()
$ws.Range("A51:H51").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new content in the order the author typed it -- this also
# drives shared-string allocation order, so keep it in this sequence.
$ws.Range("A51").Value = "No528. Random Pick with Weight"
$ws.Range("C51").Value = "https://leetcode.com/problems/random-pick-with-weight/"

# Knowledge-point tweak on row 49, made while filling in row 51.
$ws.Range("E49").Value = "数学，随机选取"

$ws.Range("B51").Value = "Medium"
$ws.Range("D51").Value = 44555
$ws.Range("E51").Value = "数学，带权重随机"
$ws.Range("F51").Value = "利用前缀和实现带有权重的随机选取，注意random取随机数的范围"
$ws.Range("G51").Value = "未复习"
$ws.Range("H51").Value = "⭕"

# Adding the hyperlink nudges the cell's style (applies a font flag), so
# restore the original copied formatting afterwards.
$ws.Hyperlinks.Add($ws.Range("C51"), "https://leetcode.com/problems/random-pick-with-weight/")
$ws.Range("A49:H49").Copy()
$ws.Range("A51:H51").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows(51).RowHeight = 42

# -- Row 52 -------------------------------------------------------------
# Row 50 (and so the new row 52) only spans columns A:G.
$ws.Range("A50:G50").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A52").Value = "81. Search in Rotated Sorted Array II"
$ws.Range("C52").Value = "https://leetcode.com/problems/search-in-rotated-sorted-array-ii/"
$ws.Range("B52").Value = "Medium"
$ws.Range("D52").Value = 44472
$ws.Range("E52").Value = "二分法"
$ws.Range("F52").Value = "当数组中存在两个递增序列时，怎么判断mid位于哪个递增序列中"
$ws.Range("G52").Value = 44555

$ws.Hyperlinks.Add($ws.Range("C52"), "https://leetcode.com/problems/search-in-rotated-sorted-array-ii/")
$ws.Range("A50:G50").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Rows(52).RowHeight = 42

# -- Selection ------------------------------------------------------------
$ws.Range("I53").Select()
